$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 1.5
$ws.Range("D2").Value = 8.300000000000001
$ws.Range("E2").Value = 7.5
$ws.Range("F2").Value = 13.5
$ws.Range("G2").Value = 1.03
$ws.Range("H2").Value = 1.79
$ws.Range("I2").Value = 0.93
$ws.Range("J2").Value = 1.69

# Row 3
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = 3.75
$ws.Range("D3").Value = 22.8
$ws.Range("E3").Value = 19.8
$ws.Range("F3").Value = 35.7
$ws.Range("G3").Value = 2.85
$ws.Range("H3").Value = 4.84
$ws.Range("I3").Value = 2.47
$ws.Range("J3").Value = 4.46

# Row 4
$ws.Range("B4").Value = 14
$ws.Range("C4").Value = 1.75
$ws.Range("D4").Value = 13.6
$ws.Range("E4").Value = 12.8
$ws.Range("F4").Value = 24.1
$ws.Range("G4").Value = 1.69
$ws.Range("H4").Value = 3.11
$ws.Range("I4").Value = 1.6
$ws.Range("J4").Value = 3.01

# Row 5
$ws.Range("B5").Value = 21
$ws.Range("C5").Value = 2.62
$ws.Range("D5").Value = 13.1
$ws.Range("E5").Value = 12.3
$ws.Range("F5").Value = 22.6
$ws.Range("G5").Value = 1.64
$ws.Range("H5").Value = 2.92
$ws.Range("I5").Value = 1.54
$ws.Range("J5").Value = 2.82

# Row 6
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = 1.25
$ws.Range("D6").Value = 10.9
$ws.Range("E6").Value = 8.5
$ws.Range("F6").Value = 15.6
$ws.Range("G6").Value = 1.36
$ws.Range("H6").Value = 2.25
$ws.Range("I6").Value = 1.07
$ws.Range("J6").Value = 1.96

# Row 7
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 0.75
$ws.Range("D7").Value = 8.9
$ws.Range("E7").Value = 8.1
$ws.Range("F7").Value = 14.5
$ws.Range("G7").Value = 1.11
$ws.Range("H7").Value = 1.91
$ws.Range("I7").Value = 1.02
$ws.Range("J7").Value = 1.81

# Row 8
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 0.87
$ws.Range("D8").Value = 10.7
$ws.Range("E8").Value = 9.9
$ws.Range("F8").Value = 17.8
$ws.Range("G8").Value = 1.34
$ws.Range("H8").Value = 2.32
$ws.Range("I8").Value = 1.24
$ws.Range("J8").Value = 2.22

# Row 9
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = 0.87
$ws.Range("D9").Value = 11.1
$ws.Range("E9").Value = 11.1
$ws.Range("F9").Value = 20.5
$ws.Range("G9").Value = 1.39
$ws.Range("H9").Value = 2.56
$ws.Range("I9").Value = 1.39
$ws.Range("J9").Value = 2.56

# Row 10
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 1.87
$ws.Range("D10").Value = 13.6
$ws.Range("E10").Value = 12
$ws.Range("F10").Value = 22.6
$ws.Range("G10").Value = 1.7
$ws.Range("H10").Value = 3.02
$ws.Range("I10").Value = 1.51
$ws.Range("J10").Value = 2.83

# Row 11
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 1.5
$ws.Range("D11").Value = 10.1
$ws.Range("E11").Value = 10.1
$ws.Range("F11").Value = 18.4
$ws.Range("G11").Value = 1.27
$ws.Range("H11").Value = 2.3
$ws.Range("I11").Value = 1.27
$ws.Range("J11").Value = 2.3

# Row 12
$ws.Range("B12").Value = 18
$ws.Range("C12").Value = 2.25
$ws.Range("D12").Value = 15.1
$ws.Range("E12").Value = 12.8
$ws.Range("F12").Value = 21.5
$ws.Range("G12").Value = 1.89
$ws.Range("H12").Value = 2.98
$ws.Range("I12").Value = 1.6
$ws.Range("J12").Value = 2.69

# Row 13
$ws.Range("B13").Value = 9
$ws.Range("C13").Value = 1.13
$ws.Range("D13").Value = 11.1
$ws.Range("E13").Value = 8.800000000000001
$ws.Range("F13").Value = 15.4
$ws.Range("G13").Value = 1.39
$ws.Range("H13").Value = 2.21
$ws.Range("I13").Value = 1.1
$ws.Range("J13").Value = 1.92

# Row 14
$ws.Range("B14").Value = 16
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 15.1
$ws.Range("E14").Value = 14.3
$ws.Range("F14").Value = 23.9
$ws.Range("G14").Value = 1.89
$ws.Range("H14").Value = 3.09
$ws.Range("I14").Value = 1.79
$ws.Range("J14").Value = 2.99

# Row 15
$ws.Range("B15").Value = 8
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 8.300000000000001
$ws.Range("E15").Value = 6.7
$ws.Range("F15").Value = 11.8
$ws.Range("G15").Value = 1.04
$ws.Range("H15").Value = 1.68
$ws.Range("I15").Value = 0.84
$ws.Range("J15").Value = 1.47

# Row 16
$ws.Range("B16").Value = 13
$ws.Range("C16").Value = 1.62
$ws.Range("D16").Value = 11.6
$ws.Range("E16").Value = 10.1
$ws.Range("F16").Value = 18.6
$ws.Range("G16").Value = 1.45
$ws.Range("H16").Value = 2.53
$ws.Range("I16").Value = 1.26
$ws.Range("J16").Value = 2.33

# Row 17
$ws.Range("B17").Value = 11
$ws.Range("C17").Value = 1.37
$ws.Range("D17").Value = 9.699999999999999
$ws.Range("E17").Value = 8.9
$ws.Range("F17").Value = 14.8
$ws.Range("G17").Value = 1.21
$ws.Range("H17").Value = 1.95
$ws.Range("I17").Value = 1.11
$ws.Range("J17").Value = 1.85

# Row 18
$ws.Range("B18").Value = 12
$ws.Range("C18").Value = 1.5
$ws.Range("D18").Value = 11.1
$ws.Range("E18").Value = 8.699999999999999
$ws.Range("F18").Value = 15.7
$ws.Range("G18").Value = 1.38
$ws.Range("H18").Value = 2.25
$ws.Range("I18").Value = 1.09
$ws.Range("J18").Value = 1.96

# Row 19
$ws.Range("B19").Value = 9
$ws.Range("C19").Value = 1.13
$ws.Range("D19").Value = 10.3
$ws.Range("E19").Value = 9.5
$ws.Range("F19").Value = 17.7
$ws.Range("G19").Value = 1.29
$ws.Range("H19").Value = 2.31
$ws.Range("I19").Value = 1.19
$ws.Range("J19").Value = 2.21

